$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.567.71"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "'2.383.85"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'507.23"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'131.39"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.547"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "'2.397.35"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").Value = "'0.0992"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "'4.86"
$ws.Range("E12").Value = "  +5.82%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "'2.809.02"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'56.522.08"
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "'21.71"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "'2.396.15"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").Value = "'10.15"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'4.06"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "'310.98"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'6.32"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'66.32"
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("D25").Value = "'0.996"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'0.373"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("E27").Value = "  -4.15%  "
$ws.Range("D28").Value = "'7.31"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").Value = "'172.91"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").Value = "'0.0₃0720"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'17.76"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").Value = "'0.834"
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("D40").Value = "'36.60"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("D42").Value = "'3.41"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("D43").Value = "'130.07"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("D44").Value = "'4.97"
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("D45").Value = "'0.570"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("D46").Value = "'0.0901"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "'242.13"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'0.0209"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "'17.17"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("E51").Value = "  -1.37%  "
